$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.820.88"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.661.36"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "242.41"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +20.91%  "
$ws.Range("D7").Value = "656.16"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.659.68"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("D12").Value = "44.25"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "0.205"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "6.51"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "4.335.34"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "96.614.06"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "3.647.18"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "8.15"
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").Value = "18.35"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("D22").Value = "0.532"
$ws.Range("E22").Value = "  +7.68%  "
$ws.Range("D23").Value = "512.27"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "6.87"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "101.42"
$ws.Range("D28").Value = "13.04"
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  +13.38%  "
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Value = "11.87"
$ws.Range("E31").Value = "  +3.36%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "33.14"
$ws.Range("E34").Value = "  +4.72%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "1.75"
$ws.Range("E36").Value = "  +7.15%  "
$ws.Range("E37").Value = "  +3.48%  "
$ws.Range("D38").Value = "615.43"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "8.75"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").Value = "42.13"
$ws.Range("E40").Value = "  +22.96%  "
$ws.Range("D41").Value = "0.159"
$ws.Range("E41").Value = "  +5.60%  "
$ws.Range("D42").Value = "0.953"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("D43").Value = "1.93"
$ws.Range("E43").Value = "  +5.76%  "
$ws.Range("D45").Value = "6.15"
$ws.Range("E45").Value = "  +7.54%  "
$ws.Range("D46").Value = "0.0444"
$ws.Range("E46").Value = "  +6.31%  "
$ws.Range("D47").Value = "0.421"
$ws.Range("E47").Value = "  +25.74%  "
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("D51").Value = "54.59"
$ws.Range("E51").Value = "  +2.77%  "
